# Auto-generated edit script: update crypto price (D) and volume-change (E) columns
# for rows 2-51 per the commit diff. D-column values are forced to remain
# plain text (matching the original inlineStr storage) even when the new
# value looks numeric, by temporarily applying a Text number format around
# the write and then resetting the style back to Normal so no stray
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.433.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.649.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3239"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07009"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.955"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.591"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.659.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001038"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06561"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9984"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.914"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.431.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.460"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.317"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -17.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.834.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.184"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.678"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08434"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.650"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -13.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.199"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.262"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06009"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02224"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2058"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.070"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -13.88%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5885"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.760"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5601"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.939"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06899"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.177"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.92%  "
